$wb = $excel.ActiveWorkbook

# --- Rename the two entity-specific sheets ---
$wsTypeTest = $wb.Worksheets.Item("V1_API_TypeTest")
$wsTypeTest.Name = "V1_API_TypeTestAPIV1"

$wsTypeTestRef = $wb.Worksheets.Item("V1_API_TypeTestRef")
$wsTypeTestRef.Name = "V1_API_TypeTestRefAPIV1"

# --- entities sheet: update entity "name"/"id" values ---
$wsEntities = $wb.Worksheets.Item("entities")
$wsEntities.Range("A2").Value = "TypeTestAPIV1"
$wsEntities.Range("C2").Value = "TypeTestAPIV1"
$wsEntities.Range("A3").Value = "TypeTestRefAPIV1"
$wsEntities.Range("C3").Value = "TypeTestRefAPIV1"
$wsEntities.Range("A4").Value = "LocationAPIV1"
$wsEntities.Range("C4").Value = "LocationAPIV1"
$wsEntities.Range("A5").Value = "PersonAPIV1"
$wsEntities.Range("C5").Value = "PersonAPIV1"

# --- attributes sheet: update "entity" (B) / "refEntity" (D) full names ---
$wsAttributes = $wb.Worksheets.Item("attributes")

# Rows 2-51 (except special rows) use V1_API_TypeTest -> V1_API_TypeTestAPIV1 in column B
for ($r = 4; $r -le 51; $r++) {
    $wsAttributes.Cells.Item($r, 2).Value = "V1_API_TypeTestAPIV1"
}

# Rows referencing TypeTestRef entity (B2, B3)
$wsAttributes.Cells.Item(2, 2).Value = "V1_API_TypeTestRefAPIV1"
$wsAttributes.Cells.Item(3, 2).Value = "V1_API_TypeTestRefAPIV1"

# refEntity column (D) referencing TypeTestRef
$typeTestRefDRows = @(10, 11, 12, 13, 36, 37, 42, 43, 48)
foreach ($r in $typeTestRefDRows) {
    $wsAttributes.Cells.Item($r, 4).Value = "V1_API_TypeTestRefAPIV1"
}

# refEntity column (D) / entity column (B) referencing Location
$wsAttributes.Cells.Item(50, 4).Value = "V1_API_LocationAPIV1"
$wsAttributes.Cells.Item(52, 2).Value = "V1_API_LocationAPIV1"
$wsAttributes.Cells.Item(53, 2).Value = "V1_API_LocationAPIV1"

# entity column (B) referencing Person
$wsAttributes.Cells.Item(54, 2).Value = "V1_API_PersonAPIV1"
$wsAttributes.Cells.Item(55, 2).Value = "V1_API_PersonAPIV1"
$wsAttributes.Cells.Item(56, 2).Value = "V1_API_PersonAPIV1"

# --- Make the renamed V1_API_TypeTestAPIV1 sheet the active tab ---
$wsTypeTest.Activate()

Write-Output "done"
